# AgentLoginLogoutReportData.xlsx edit
# - Update the two "Login/Logout" sample date/time strings on the
#   "Queries" sheet (8th sheet, cells D2/E2) from Feb-2020 values to
#   Mar-2020 values.
# - Move the active cell selection on that sheet from F2 to E2.
# - Update the auto-fit row height for row 2 (wrapped long SQL query
#   in F2) to the recalculated height.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)

# The cells hold text (quote-prefixed) values that look like dates, so
# prefix the new values with a leading apostrophe to keep them stored
# as literal text (preserves the existing "quote prefix" cell style).
$ws.Cells.Item(2, 4).Value2 = "'19-03-2020 00:00:00"
$ws.Cells.Item(2, 5).Value2 = "'27-03-2020 00:00:00"

# Move the selection/active cell from F2 to E2 on this (active) sheet.
[void]$ws.Activate()
[void]$ws.Range("E2").Select()

# Row 2 autofits to a slightly shorter height after the edit.
$ws.Rows.Item(2).RowHeight = 362.5
